# Config.xlsx update - "Upto Group 3 MGT 7"
# Updates file paths in the config table and appends a new block of
# Key/Value configuration rows (rows 11-29) used for MGT-7 Group 3
# (principal business activities / director shareholdings / holdings etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update the two path values (now pointing at the new repo checkout
#    location, and the mapping-config filename lost its "-1" suffix).
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "/Users/gundukalyan/Documents/GitHub/mns-json-prep/DataExtraction/Input/Form MGT-7-22092022_signed - Json data L&T.pdf"
$ws.Range("B4").Value = "/Users/gundukalyan/Documents/GitHub/mns-json-prep/DataExtraction/Input/MGT7_Newmapping_config.xlsx"

# ---------------------------------------------------------------------
# 2) Un-bold the header row (A1:E1) - the bold Arial font is dropped.
# ---------------------------------------------------------------------
$ws.Range("A1:E1").Font.Bold = $false

# ---------------------------------------------------------------------
# 3) Fill in the new Key/Value rows. Most of these cells already carry
#    the correct banded row style, so we only need to set their values.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "principal_business_activities_field_name"
$ws.Range("B11").Value = "principal_business_activities"

$ws.Range("A12").Value = "director_shareholdings_field_name"
$ws.Range("B12").Value = "director_shareholdings"

$ws.Range("A14").Value = "Hold_Sub_Assoc_field_name"
$ws.Range("B14").Value = "Holding/ Subsidiary/Associate/  Joint Venture"

$ws.Range("A15").Value = "director_remuneration_field_name"
$ws.Range("B15").Value = "director_remuneration"

$ws.Range("A17").Value = "principal_business_activities_table_name"
$ws.Range("B17").Value = "principal_business_activities"

$ws.Range("A19").Value = "year_field_name"
$ws.Range("B19").Value = "Year"

$ws.Range("A21").Value = "Hold_Sub_Assoc_column_name"
$ws.Range("B21").Value = "HOLD_SUB_ASSOC"

$ws.Range("A23").Value = "cin_column_name_in_db"
$ws.Range("B23").Value = "cin"

$ws.Range("A24").Value = "company_name_column_name_in_db"
$ws.Range("B24").Value = "company_name"

$ws.Range("A26").Value = "associate_keyword_in_xml"
$ws.Range("B26").Value = "ASSOC"

$ws.Range("A27").Value = "holding_keyword_in_xml"
$ws.Range("B27").Value = "HOLD"

$ws.Range("A28").Value = "joint_venture_keyword_in_xml"
$ws.Range("B28").Value = "JOINT"

$ws.Range("A29").Value = "subsidiary_keyword_in_xml"
$ws.Range("B29").Value = "SUBS"

# ---------------------------------------------------------------------
# 4) A handful of cells picked up a slightly different banding style
#    while the table was being extended in Excel - replicate that.
# ---------------------------------------------------------------------

# A13 switches from the "grey separator" look to the plain style used by
# A11/A15 (no value is present in row 13, only the format changed).
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# A27 / B27 switch the same way.
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B27").PasteSpecial(-4122) | Out-Null

# A29 flips to the alternate band used by A12/A28.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null

# A26 keeps the A24-style band but with the wrap turned off, bottom
# vertical alignment, and a slightly lighter grey left border - this is
# the one cell that ends up on a brand-new style in the saved file.
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A26").PasteSpecial(-4122) | Out-Null
$ws.Range("A26").WrapText = $false
$ws.Range("A26").VerticalAlignment = -4107   # xlBottom
$ws.Range("A26").Borders.Color = 13158600    # RGB(200,200,200) on every edge
$ws.Range("A26").Borders.Item(7).Color = 11184810   # xlEdgeLeft -> RGB(170,170,170)

$excel.CutCopyMode = 0
